$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.856.51"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.040.79"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'250.68"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'57.79"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'61.33"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  +3.25%  "
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").Value = "'16.25"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("D14").Value = "2.342.43"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("E15").Value = "  -5.86%  "
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("D17").Value = "2.044.17"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "36.830.46"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'16.66"
$ws.Range("E19").Value = "  +12.67%  "
$ws.Range("D20").Value = "'75.06"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("E21").Value = "  +5.84%  "
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "'236.59"
$ws.Range("E23").Value = "  -1.98%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  +15.95%  "
$ws.Range("D27").Value = "'169.01"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'9.21"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "'20.14"
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").Value = "'4.43"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -4.40%  "
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").Value = "'0.112"
$ws.Range("E39").Value = "  +12.59%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").Value = "'17.72"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").Value = "'96.52"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  +14.91%  "
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("D48").Value = "1.281.69"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'6.72"
$ws.Range("D51").Value = "2.233.03"
$ws.Range("E51").Value = "  -1.87%  "
